$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.671.09"
$ws.Range("E2").Value = "  +1.80%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.868.69"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "235.69"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "

# Row 6 - USDC
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.22%  "

# Row 8 - Cardano
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2762"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.19%  "

# Row 9 - Dogecoin
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06384"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10 - Solana
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "18.05"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +11.60%  "

# Row 11 - WrappedEther
$ws.Range("D11").Value = "1.876.48"
$ws.Range("E11").Value = "  +0.07%  "

# Row 12 - TRON
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07452"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

# Row 13 - Polkadot
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.977"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "

# Row 14 - Litecoin
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "85.29"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.45%  "

# Row 15 - Polygon
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6360"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "30.626.64"
$ws.Range("E16").Value = "  +1.84%  "

# Row 17 - BitcoinCash
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "247.12"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +7.73%  "

# Row 18 - Dai
$ws.Range("E18").Value = "  -0.15%  "

# Row 19 - Avalanche
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.80"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.62%  "

# Row 20 - ShibaInu
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007406"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "

# Row 21 - BinanceUSD
$ws.Range("E21").Value = "  -0.09%  "

# Row 22 - Uniswap
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.972"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.86%  "

# Row 23 - Chainlink
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.090"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.05%  "

# Row 24 - Cosmos
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "9.396"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "

# Row 25 - Monero
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "164.24"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "

# Row 26 - EthereumClassic
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "18.32"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.30%  "

# Row 27 - LidoDAOToken
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.909"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.78%  "

# Row 28 - Stellar
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.1022"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.49%  "

# Row 29 - Toncoin
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.382"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30 - InternetComputer(DFINITY)
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.088"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.70%  "

# Row 31 - Filecoin
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.874"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.55%  "

# Row 32 - Hedera
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.04943"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.52%  "

# Row 33 - ARBITRUM
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.159"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.26%  "

# Row 34 - ImmutableX
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7112"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "

# Row 35 - HuobiToken
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.712"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.68%  "

# Row 36 - VeChain
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.01913"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.34%  "

# Row 37 - MXToken
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.694"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.40%  "

# Row 38 - TrustWalletToken
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.8843"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.42%  "

# Row 39 - RenderToken
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.009"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.95%  "

# Row 40 - Quant
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "105.95"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.55%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.04%  "

# Row 42 - TheSandbox
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.4111"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.96%  "

# Row 43 - FraxShare
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.556"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "

# Row 44 - Aptos
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "7.351"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.94%  "

# Row 45 - Aave
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "62.52"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.20%  "

# Row 46 - Algorand
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.1226"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.75%  "

# Row 47 - EnergySwap
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.688"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "

# Row 48 - Elrond
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "33.73"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.97%  "

# Row 49 - Cronos (was NEARProtocol)
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.05590"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.40%  "

# Row 50 - NEARProtocol (was Cronos)
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.386"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.30%  "

# Row 51 - Decentraland
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.3714"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.94%  "

